$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the "Förändrad" (changed-on) column C for rows 2-28: 45426 -> 45427
for ($r = 2; $r -le 28; $r++) {
    $ws.Cells.Item($r, 3).Value = 45427
    $ws.Cells.Item($r, 3).NumberFormat = "YYYY-MM-DD"
}

# --- 2. Row 28 gains an explicit row height (ht="15" customHeight="1")
$ws.Rows.Item(28).RowHeight = 15

# --- 3. Add new row 29: A 18860-2024
$ws.Rows.Item(29).RowHeight = 15
$ws.Cells.Item(29, 1).Value = "A 18860-2024"
$ws.Cells.Item(29, 2).Value = 45426
$ws.Cells.Item(29, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(29, 3).Value = 45427
$ws.Cells.Item(29, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(29, 4).Value = "OKÄNT"
$ws.Cells.Item(29, 5).Value = "OKÄNT"
$ws.Cells.Item(29, 6).Value = "Sveaskog"
$ws.Cells.Item(29, 7).Value = 9.6
$ws.Cells.Item(29, 8).Value = 0
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 10).Value = 0
$ws.Cells.Item(29, 11).Value = 0
$ws.Cells.Item(29, 12).Value = 0
$ws.Cells.Item(29, 13).Value = 0
$ws.Cells.Item(29, 14).Value = 0
$ws.Cells.Item(29, 15).Value = 0
$ws.Cells.Item(29, 16).Value = 0
$ws.Cells.Item(29, 17).Value = 0
$ws.Cells.Item(29, 18).WrapText = $true

# --- 4. Add new row 30: A 18855-2024 (no row height override, no F value)
$ws.Cells.Item(30, 1).Value = "A 18855-2024"
$ws.Cells.Item(30, 2).Value = 45426
$ws.Cells.Item(30, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(30, 3).Value = 45427
$ws.Cells.Item(30, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(30, 4).Value = "OKÄNT"
$ws.Cells.Item(30, 5).Value = "OKÄNT"
$ws.Cells.Item(30, 7).Value = 5
$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(30, 9).Value = 0
$ws.Cells.Item(30, 10).Value = 0
$ws.Cells.Item(30, 11).Value = 0
$ws.Cells.Item(30, 12).Value = 0
$ws.Cells.Item(30, 13).Value = 0
$ws.Cells.Item(30, 14).Value = 0
$ws.Cells.Item(30, 15).Value = 0
$ws.Cells.Item(30, 16).Value = 0
$ws.Cells.Item(30, 17).Value = 0
$ws.Cells.Item(30, 18).WrapText = $true
